$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 17:48:33"
$ws.Cells.Item(3,1).Value = "Total filas: 291"
$ws.Cells.Item(40,1).Value = "06:35:22"
$ws.Cells.Item(40,2).Value = "08:29"
$ws.Cells.Item(40,3).Value = "15_ABASTO"
$ws.Cells.Item(40,4).Value = 114
$ws.Cells.Item(40,5).Value = "LP1912"
$ws.Cells.Item(41,1).Value = "06:35:22"
$ws.Cells.Item(41,2).Value = "08:29"
$ws.Cells.Item(41,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(41,4).Value = 114
$ws.Cells.Item(41,5).Value = "LP1912"
$ws.Cells.Item(58,1).Value = "08:10:18"
$ws.Cells.Item(58,2).Value = "09:18"
$ws.Cells.Item(58,3).Value = "14_ABASTO"
$ws.Cells.Item(58,4).Value = 68
$ws.Cells.Item(58,5).Value = "LP1912"
$ws.Cells.Item(59,1).Value = "07:38:39"
$ws.Cells.Item(59,2).Value = "09:18"
$ws.Cells.Item(59,3).Value = "15X38_ABASTO"
$ws.Cells.Item(59,4).Value = 100
$ws.Cells.Item(59,5).Value = "LP1912"
$ws.Cells.Item(77,1).Value = "09:23:23"
$ws.Cells.Item(77,2).Value = "10:29"
$ws.Cells.Item(77,3).Value = "14_ABASTO"
$ws.Cells.Item(77,4).Value = 66
$ws.Cells.Item(77,5).Value = "LP1912"
$ws.Cells.Item(78,1).Value = "08:37:25"
$ws.Cells.Item(78,2).Value = "10:29"
$ws.Cells.Item(78,3).Value = "15_ABASTO"
$ws.Cells.Item(78,4).Value = 112
$ws.Cells.Item(78,5).Value = "LP1912"
$ws.Cells.Item(86,1).Value = "10:50:41"
$ws.Cells.Item(86,2).Value = "10:59"
$ws.Cells.Item(86,3).Value = "10_OLMOS"
$ws.Cells.Item(86,4).Value = 9
$ws.Cells.Item(86,5).Value = "LP1912"
$ws.Cells.Item(87,1).Value = "09:23:23"
$ws.Cells.Item(87,2).Value = "10:59"
$ws.Cells.Item(87,3).Value = "27_EL RETIRO"
$ws.Cells.Item(87,4).Value = 96
$ws.Cells.Item(87,5).Value = "LP1912"
$ws.Cells.Item(130,1).Value = "11:47:17"
$ws.Cells.Item(130,2).Value = "12:33"
$ws.Cells.Item(130,3).Value = "14_ABASTO"
$ws.Cells.Item(130,4).Value = 46
$ws.Cells.Item(130,5).Value = "LP1912"
$ws.Cells.Item(131,1).Value = "11:34:59"
$ws.Cells.Item(131,2).Value = "12:33"
$ws.Cells.Item(131,3).Value = "15_ABASTO"
$ws.Cells.Item(131,4).Value = 59
$ws.Cells.Item(131,5).Value = "LP1912"
$ws.Cells.Item(146,1).Value = "11:47:17"
$ws.Cells.Item(146,2).Value = "13:03"
$ws.Cells.Item(146,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(146,4).Value = 76
$ws.Cells.Item(146,5).Value = "LP1912"
$ws.Cells.Item(147,1).Value = "11:34:59"
$ws.Cells.Item(147,2).Value = "13:03"
$ws.Cells.Item(147,3).Value = "215C_EL PATO"
$ws.Cells.Item(147,4).Value = 89
$ws.Cells.Item(147,5).Value = "LP1912"
$ws.Cells.Item(160,1).Value = "12:11:52"
$ws.Cells.Item(160,2).Value = "13:32"
$ws.Cells.Item(160,3).Value = "14_ABASTO"
$ws.Cells.Item(160,4).Value = 81
$ws.Cells.Item(160,5).Value = "LP1912"
$ws.Cells.Item(161,1).Value = "11:34:59"
$ws.Cells.Item(161,2).Value = "13:32"
$ws.Cells.Item(161,3).Value = "215A_EL PATO"
$ws.Cells.Item(161,4).Value = 118
$ws.Cells.Item(161,5).Value = "LP1912"
$ws.Cells.Item(169,1).Value = "13:14:29"
$ws.Cells.Item(169,2).Value = "14:02"
$ws.Cells.Item(169,3).Value = "16_SANTA ANA"
$ws.Cells.Item(169,4).Value = 48
$ws.Cells.Item(169,5).Value = "LP1912"
$ws.Cells.Item(170,1).Value = "12:33:21"
$ws.Cells.Item(170,2).Value = "14:02"
$ws.Cells.Item(170,3).Value = "10_OLMOS"
$ws.Cells.Item(170,4).Value = 89
$ws.Cells.Item(170,5).Value = "LP1912"
$ws.Cells.Item(260,1).Value = "17:48:33"
$ws.Cells.Item(260,2).Value = "17:48"
$ws.Cells.Item(260,3).Value = "215B_EL PATO"
$ws.Cells.Item(260,4).Value = 0
$ws.Cells.Item(260,5).Value = "LP1912"
$ws.Cells.Item(261,1).Value = "17:48:33"
$ws.Cells.Item(261,2).Value = "17:48"
$ws.Cells.Item(261,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(261,4).Value = 0
$ws.Cells.Item(261,5).Value = "LP1912"
$ws.Cells.Item(262,1).Value = "16:28:03"
$ws.Cells.Item(262,2).Value = "17:49"
$ws.Cells.Item(262,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(262,4).Value = 81
$ws.Cells.Item(262,5).Value = "LP1912"
$ws.Cells.Item(263,1).Value = "16:37:06"
$ws.Cells.Item(263,2).Value = "17:50"
$ws.Cells.Item(263,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(263,4).Value = 73
$ws.Cells.Item(263,5).Value = "LP1912"
$ws.Cells.Item(264,1).Value = "17:13:39"
$ws.Cells.Item(264,2).Value = "17:51"
$ws.Cells.Item(264,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(264,4).Value = 38
$ws.Cells.Item(264,5).Value = "LP1912"
$ws.Cells.Item(265,1).Value = "16:37:06"
$ws.Cells.Item(265,2).Value = "17:52"
$ws.Cells.Item(265,3).Value = "10_OLMOS"
$ws.Cells.Item(265,4).Value = 75
$ws.Cells.Item(265,5).Value = "LP1912"
$ws.Cells.Item(266,1).Value = "16:28:03"
$ws.Cells.Item(266,2).Value = "17:53"
$ws.Cells.Item(266,3).Value = "10_OLMOS"
$ws.Cells.Item(266,4).Value = 85
$ws.Cells.Item(266,5).Value = "LP1912"
$ws.Cells.Item(267,1).Value = "16:28:03"
$ws.Cells.Item(267,2).Value = "17:58"
$ws.Cells.Item(267,3).Value = "17_ROMERO"
$ws.Cells.Item(267,4).Value = 90
$ws.Cells.Item(267,5).Value = "LP1912"
$ws.Cells.Item(268,1).Value = "16:14:44"
$ws.Cells.Item(268,2).Value = "18:00"
$ws.Cells.Item(268,3).Value = "10_OLMOS"
$ws.Cells.Item(268,4).Value = 106
$ws.Cells.Item(268,5).Value = "LP1912"
$ws.Cells.Item(269,1).Value = "16:14:44"
$ws.Cells.Item(269,2).Value = "18:05"
$ws.Cells.Item(269,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(269,4).Value = 111
$ws.Cells.Item(269,5).Value = "LP1912"
$ws.Cells.Item(270,1).Value = "16:28:03"
$ws.Cells.Item(270,2).Value = "18:06"
$ws.Cells.Item(270,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(270,4).Value = 98
$ws.Cells.Item(270,5).Value = "LP1912"
$ws.Cells.Item(271,1).Value = "16:14:44"
$ws.Cells.Item(271,2).Value = "18:10"
$ws.Cells.Item(271,3).Value = "15_ABASTO"
$ws.Cells.Item(271,4).Value = 116
$ws.Cells.Item(271,5).Value = "LP1912"
$ws.Cells.Item(272,1).Value = "16:14:44"
$ws.Cells.Item(272,2).Value = "18:10"
$ws.Cells.Item(272,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(272,4).Value = 116
$ws.Cells.Item(272,5).Value = "LP1912"
$ws.Cells.Item(273,1).Value = "16:28:03"
$ws.Cells.Item(273,2).Value = "18:17"
$ws.Cells.Item(273,3).Value = "10_OLMOS"
$ws.Cells.Item(273,4).Value = 109
$ws.Cells.Item(273,5).Value = "LP1912"
$ws.Cells.Item(274,1).Value = "16:52:42"
$ws.Cells.Item(274,2).Value = "18:21"
$ws.Cells.Item(274,3).Value = "215C_EL PATO"
$ws.Cells.Item(274,4).Value = 89
$ws.Cells.Item(274,5).Value = "LP1912"
$ws.Cells.Item(275,1).Value = "16:28:03"
$ws.Cells.Item(275,2).Value = "18:22"
$ws.Cells.Item(275,3).Value = "215C_EL PATO"
$ws.Cells.Item(275,4).Value = 114
$ws.Cells.Item(275,5).Value = "LP1912"
$ws.Cells.Item(276,1).Value = "16:28:03"
$ws.Cells.Item(276,2).Value = "18:25"
$ws.Cells.Item(276,3).Value = "16_SANTA ANA"
$ws.Cells.Item(276,4).Value = 117
$ws.Cells.Item(276,5).Value = "LP1912"
$ws.Cells.Item(277,1).Value = "16:52:42"
$ws.Cells.Item(277,2).Value = "18:29"
$ws.Cells.Item(277,3).Value = "14_ABASTO"
$ws.Cells.Item(277,4).Value = 97
$ws.Cells.Item(277,5).Value = "LP1912"
$ws.Cells.Item(278,1).Value = "16:37:06"
$ws.Cells.Item(278,2).Value = "18:30"
$ws.Cells.Item(278,3).Value = "14_ABASTO"
$ws.Cells.Item(278,4).Value = 113
$ws.Cells.Item(278,5).Value = "LP1912"
$ws.Cells.Item(279,1).Value = "17:48:33"
$ws.Cells.Item(279,2).Value = "18:32"
$ws.Cells.Item(279,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(279,4).Value = 44
$ws.Cells.Item(279,5).Value = "LP1912"
$ws.Cells.Item(280,1).Value = "16:37:06"
$ws.Cells.Item(280,2).Value = "18:36"
$ws.Cells.Item(280,3).Value = "15X38_ABASTO"
$ws.Cells.Item(280,4).Value = 119
$ws.Cells.Item(280,5).Value = "LP1912"
$ws.Cells.Item(281,1).Value = "17:13:39"
$ws.Cells.Item(281,2).Value = "18:36"
$ws.Cells.Item(281,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(281,4).Value = 83
$ws.Cells.Item(281,5).Value = "LP1912"
$ws.Cells.Item(282,1).Value = "17:36:10"
$ws.Cells.Item(282,2).Value = "18:37"
$ws.Cells.Item(282,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(282,4).Value = 61
$ws.Cells.Item(282,5).Value = "LP1912"
$ws.Cells.Item(283,1).Value = "17:13:39"
$ws.Cells.Item(283,2).Value = "18:41"
$ws.Cells.Item(283,3).Value = "10_OLMOS"
$ws.Cells.Item(283,4).Value = 88
$ws.Cells.Item(283,5).Value = "LP1912"
$ws.Cells.Item(284,1).Value = "16:52:42"
$ws.Cells.Item(284,2).Value = "18:45"
$ws.Cells.Item(284,3).Value = "16_SANTA ANA"
$ws.Cells.Item(284,4).Value = 113
$ws.Cells.Item(284,5).Value = "LP1912"
$ws.Cells.Item(285,1).Value = "17:13:39"
$ws.Cells.Item(285,2).Value = "18:52"
$ws.Cells.Item(285,3).Value = "17_ROMERO"
$ws.Cells.Item(285,4).Value = 99
$ws.Cells.Item(285,5).Value = "LP1912"
$ws.Cells.Item(286,1).Value = "17:13:39"
$ws.Cells.Item(286,2).Value = "18:57"
$ws.Cells.Item(286,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(286,4).Value = 104
$ws.Cells.Item(286,5).Value = "LP1912"
$ws.Cells.Item(287,1).Value = "17:13:39"
$ws.Cells.Item(287,2).Value = "18:59"
$ws.Cells.Item(287,3).Value = "14_ABASTO"
$ws.Cells.Item(287,4).Value = 106
$ws.Cells.Item(287,5).Value = "LP1912"
$ws.Cells.Item(288,1).Value = "17:36:10"
$ws.Cells.Item(288,2).Value = "19:00"
$ws.Cells.Item(288,3).Value = "14_ABASTO"
$ws.Cells.Item(288,4).Value = 84
$ws.Cells.Item(288,5).Value = "LP1912"
$ws.Cells.Item(289,1).Value = "17:13:39"
$ws.Cells.Item(289,2).Value = "19:03"
$ws.Cells.Item(289,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(289,4).Value = 110
$ws.Cells.Item(289,5).Value = "LP1912"
$ws.Cells.Item(290,1).Value = "17:36:10"
$ws.Cells.Item(290,2).Value = "19:04"
$ws.Cells.Item(290,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(290,4).Value = 88
$ws.Cells.Item(290,5).Value = "LP1912"
$ws.Cells.Item(291,1).Value = "17:48:33"
$ws.Cells.Item(291,2).Value = "19:12"
$ws.Cells.Item(291,3).Value = "27_EL RETIRO"
$ws.Cells.Item(291,4).Value = 84
$ws.Cells.Item(291,5).Value = "LP1912"
$ws.Cells.Item(292,1).Value = "17:36:10"
$ws.Cells.Item(292,2).Value = "19:17"
$ws.Cells.Item(292,3).Value = "27_EL RETIRO"
$ws.Cells.Item(292,4).Value = 101
$ws.Cells.Item(292,5).Value = "LP1912"
$ws.Cells.Item(293,1).Value = "17:36:10"
$ws.Cells.Item(293,2).Value = "19:17"
$ws.Cells.Item(293,3).Value = "14X44_ABASTO"
$ws.Cells.Item(293,4).Value = 101
$ws.Cells.Item(293,5).Value = "LP1912"
$ws.Cells.Item(294,1).Value = "17:36:10"
$ws.Cells.Item(294,2).Value = "19:28"
$ws.Cells.Item(294,3).Value = "215C_EL PATO"
$ws.Cells.Item(294,4).Value = 112
$ws.Cells.Item(294,5).Value = "LP1912"
$ws.Cells.Item(295,1).Value = "17:48:33"
$ws.Cells.Item(295,2).Value = "19:35"
$ws.Cells.Item(295,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(295,4).Value = 107
$ws.Cells.Item(295,5).Value = "LP1912"
$ws.Cells.Item(296,1).Value = "17:48:33"
$ws.Cells.Item(296,2).Value = "19:39"
$ws.Cells.Item(296,3).Value = "15X38_ABASTO"
$ws.Cells.Item(296,4).Value = 111
$ws.Cells.Item(296,5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 17:48:33"
$ws.Cells.Item(3,1).Value = "Total filas: 48"
$ws.Cells.Item(47,1).Value = "17:48:33"
$ws.Cells.Item(47,2).Value = "17:48"
$ws.Cells.Item(47,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,5).Value = "LP1912"
$ws.Cells.Item(48,1).Value = "17:48:33"
$ws.Cells.Item(48,2).Value = "17:48"
$ws.Cells.Item(48,3).Value = "215B_EL PATO"
$ws.Cells.Item(48,4).Value = 0
$ws.Cells.Item(48,5).Value = "LP1912"
$ws.Cells.Item(49,1).Value = "16:52:42"
$ws.Cells.Item(49,2).Value = "18:21"
$ws.Cells.Item(49,3).Value = "215C_EL PATO"
$ws.Cells.Item(49,4).Value = 89
$ws.Cells.Item(49,5).Value = "LP1912"
$ws.Cells.Item(50,1).Value = "16:28:03"
$ws.Cells.Item(50,2).Value = "18:22"
$ws.Cells.Item(50,3).Value = "215C_EL PATO"
$ws.Cells.Item(50,4).Value = 114
$ws.Cells.Item(50,5).Value = "LP1912"
$ws.Cells.Item(51,1).Value = "17:13:39"
$ws.Cells.Item(51,2).Value = "19:03"
$ws.Cells.Item(51,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(51,4).Value = 110
$ws.Cells.Item(51,5).Value = "LP1912"
$ws.Cells.Item(52,1).Value = "17:36:10"
$ws.Cells.Item(52,2).Value = "19:04"
$ws.Cells.Item(52,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(52,4).Value = 88
$ws.Cells.Item(52,5).Value = "LP1912"
$ws.Cells.Item(53,1).Value = "17:36:10"
$ws.Cells.Item(53,2).Value = "19:28"
$ws.Cells.Item(53,3).Value = "215C_EL PATO"
$ws.Cells.Item(53,4).Value = 112
$ws.Cells.Item(53,5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 17:48:33"
